$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("papers")
$ws.Range("A1").Value = "TEST"
